$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new worksheet "table with merged dividers" at the end (after
#    "broken stats"). It mirrors sheet4 ("table with dividers") but the
#    group-divider rows are merged across A:C (instead of being three
#    separately-styled cells), which is what exposed the divider-locating
#    bug referenced in the commit message.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add($null, $lastSheet)
$ws8.Name = "table with merged dividers"

# Column widths (~A=14.5, B=13.3, C=13 characters)
$ws8.Columns.Item(1).ColumnWidth = 13.7
$ws8.Columns.Item(2).ColumnWidth = 12.42
$ws8.Columns.Item(3).ColumnWidth = 12.2

# Header row
$ws8.Range("A1").Value = "team"
$ws8.Range("B1").Value = "plays"
$ws8.Range("C1").Value = "points"
$ws8.Range("A1:C1").Font.Bold = $true

# --- Group A -----------------------------------------------------------
$ws8.Range("A2").Value = "Group A"
$ws8.Range("A2:C3").Font.Bold = $true
$ws8.Range("A2:C3").Font.Color = 255
$ws8.Range("A2:C3").HorizontalAlignment = -4108
$ws8.Range("A2:C3").VerticalAlignment = -4108
$ws8.Range("A2:C3").Merge()

$ws8.Range("A4").Value = "PSG"
$ws8.Range("B4").Value = 3
$ws8.Range("C4").Value = 7

$ws8.Range("A5").Value = "Manchester City"
$ws8.Range("B5").Value = 3
$ws8.Range("C5").Value = 6

$ws8.Range("A6").Value = "Club Brugge"
$ws8.Range("B6").Value = 3
$ws8.Range("C6").Value = 4

$ws8.Range("A7").Value = "RB Leipzig"
$ws8.Range("B7").Value = 3
$ws8.Range("C7").Value = 0

# --- Group B -----------------------------------------------------------
$ws8.Range("A8").Value = "Group B"
$ws8.Range("A8:C9").Font.Bold = $true
$ws8.Range("A8:C9").Font.Color = 255
$ws8.Range("A8:C9").HorizontalAlignment = -4108
$ws8.Range("A8:C9").VerticalAlignment = -4108
$ws8.Range("A8:C9").Merge()

# Marker cell just past the merged divider - this is the "min/max column
# index" case the bug fix needed to handle correctly.
$ws8.Range("D8").Value = "``"

$ws8.Range("A10").Value = "Liverpool"
$ws8.Range("B10").Value = 3
$ws8.Range("C10").Value = 9

$ws8.Range("A11").Value = "Atletico Madrid"
$ws8.Range("B11").Value = 3
$ws8.Range("C11").Value = 4

$ws8.Range("A12").Value = "FC Porto"
$ws8.Range("B12").Value = 3
$ws8.Range("C12").Value = 4

$ws8.Range("A13").Value = "AC Milan"
$ws8.Range("B13").Value = 3
$ws8.Range("C13").Value = 0

$ws8.Range("D9").Select()

# ---------------------------------------------------------------------------
# 2. Re-activate "table with dividers" (sheet 4) as the active tab, with the
#    cursor left on I7 (the cell the author was inspecting while fixing the
#    divider-locating bug).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Range("I7").Select()
